# Update the "State" column (D3:D23) from "Draft ready" to "Sent"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3:D23").Value = "Sent"

# Reflect the saved cursor/selection position from the diff
$ws.Range("D24").Select()
